# Add the two new single-atom energy entries (Ag, Cu) below the existing
# Pd/Pt/Au rows, mirroring the existing pattern: element symbol in column A,
# a formula converting Hartree energy to the Ha/27.2114 value in column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value2 = "Ag"
$ws.Range("B5").Formula = "=-1004.93065412009/27.2114"

$ws.Range("A6").Value2 = "Cu"
$ws.Range("B6").Formula = "=-1305.9228912226/27.2114"

# Match the author's final selection (cell B6) recorded in the saved file.
$ws.Range("B6").Select()
